# Remove two snowballing entries that are no longer needed, preparing the
# sheet for appraisal and extraction. Deleting the entire rows shifts every
# row below upward and Excel automatically repacks the shared string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 132: "N Sethi, S Rani, P Singh" / "Ants optimization for minimal test
# case selection and prioritization..." entry.
$ws.Rows.Item(132).Delete() | Out-Null

# Row 169 in the original numbering is now row 168 after the deletion above:
# "L Mei, Y Cai, C Jia, B Jiang, WK Chan" / "Test pair selection for test
# case prioritization..." entry.
$ws.Rows.Item(168).Delete() | Out-Null

# Restore the active selection to where the author left off working
# (row 168 in the final, renumbered sheet).
$ws.Rows.Item(168).Select() | Out-Null

# Best-effort: scroll the window so row 193 is near the top, matching the
# author's last recorded view position.
$excel.ActiveWindow.ScrollRow = 193
$excel.ActiveWindow.ScrollColumn = 1
